# Add two new TODO sections ("7. Pages Controller" and "8. FrontPage")
# right after the existing "-create modal and bind variable to it" line,
# and move the "_GoBack" bookmark so it still marks the point of the most
# recently typed text (now at the end of the new "8. FrontPage" section).

$d = $word.ActiveDocument

# Locate the end of the text "-create modal and bind variable to it"
# (this is also the exact spot, just before the _GoBack bookmark, where
# the new paragraphs need to be inserted).
$rng = $d.Content
$rng.Find.Execute("bind variable to it", $false, $false, $false, $false, `
                   $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$insertStart = $rng.Start

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Build the eight new paragraphs as raw WordprocessingML so tabs come
# through as real <w:tab/> runs (matching the rest of the document) rather
# than literal tab characters. The final paragraph temporarily ends with a
# two-character "ZZ" placeholder; it lets us anchor the relocated bookmark
# just before it without hitting an edge case where bookmarking the very
# last character of the whole document snaps back to position 0. The
# placeholder is stripped immediately afterwards.
$xml = @"
<w:p $ns><w:r><w:t xml:space="preserve">7. Pages Controller: </w:t></w:r></w:p>
<w:p $ns><w:r><w:tab/><w:t>Create modal for article creation/modification</w:t></w:r></w:p>
<w:p $ns><w:r><w:tab/><w:t xml:space="preserve">Create modal for delete </w:t></w:r></w:p>
<w:p $ns><w:r><w:tab/><w:t xml:space="preserve">Bind modal for each case to a variable </w:t></w:r></w:p>
<w:p $ns><w:r><w:tab/><w:t xml:space="preserve">Change the value of the variable when a button clicked for the specific action </w:t></w:r></w:p>
<w:p $ns><w:r><w:tab/><w:t>Bind the variables in controller to modals then reset them after validation, submission of action</w:t></w:r></w:p>
<w:p $ns><w:r><w:t>8. FrontPage:</w:t></w:r></w:p>
<w:p $ns><w:r><w:tab/><w:t>Like single page application the frontend is one page that calls specific components when needed and call specific data from databaseZZ</w:t></w:r></w:p>
"@

$insertRange = $d.Range($insertStart, $insertStart)
$insertRange.InsertXML($xml)

# Re-anchor the "_GoBack" bookmark at the end of the newly typed content
# (right after "...call specific data from database", before the "ZZ"
# placeholder) instead of its old spot after "...bind variable to it".
$bmEnd = $d.Content.End - 3
$bmRange = $d.Range($bmEnd, $bmEnd)
$d.Bookmarks("_GoBack").Delete()
$d.Bookmarks.Add("_GoBack", $bmRange)

# Remove the temporary "ZZ" placeholder now that the bookmark is anchored.
$cleanupRange = $d.Range($d.Content.End - 3, $d.Content.End - 1)
$cleanupRange.Delete()
